$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newItems = @(
    "hidden objects should not be included in shade analysis?",
    "remove checkboxes to show/hide objects",
    "zoom to fit, 'F' key",
    "clear properties on de-select object",
    "object list sorted by group",
    "active area polygon rotation",
    "composite objects: regular fixed array"
)

$startRow = 67
for ($i = 0; $i -lt $newItems.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = "Not done"
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Value = $newItems[$i]
    $cellB.Interior.Color = 11916796
}

$ws.Range("F79").Select()
